# Insert a new weekly record at row 139 of the Zapallo italiano sheet.
# All existing rows from 139 downward shift down by one row (139->140, ..., 215->216),
# and the sheet dimension grows from A1:R215 to A1:R216.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 139, pushing old rows 139..215 to 140..216
$ws.Rows.Item(139).Insert()

# Populate the new row 139 with the new weekly data point
$ws.Cells.Item(139, 1).Value = 11
$ws.Cells.Item(139, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(139, 3).Value = "Bíobío"
$ws.Cells.Item(139, 4).Value = 45029
$ws.Cells.Item(139, 5).Value = 8
$ws.Cells.Item(139, 6).Value = 100112032
$ws.Cells.Item(139, 7).Value = "Zapallo italiano"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 220
$ws.Cells.Item(139, 11).Value = 6500
$ws.Cells.Item(139, 12).Value = 7000
$ws.Cells.Item(139, 13).Value = 6773
$ws.Cells.Item(139, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(139, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(139, 16).Value = 135
$ws.Cells.Item(139, 17).Value = 50
$ws.Cells.Item(139, 18).Value = "Hortaliza"

# Make sure the new row's date cell carries the same date style/format as the rest of column D
$ws.Cells.Item(139, 4).NumberFormat = $ws.Cells.Item(140, 4).NumberFormat
